$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Name = "Main_240603"

$ws.Range("F3").Value = "등산용 대형 배낭입니다. 산업계 군대의 보병에게 지급되는 종류의 것과 기능상으로 동일합니다."
$ws.Range("F2").Value = "대형 배낭"

$range = $ws.Range("F2:F3")
$range.FormatConditions.Delete()
$fc = $range.FormatConditions.Add(2, 3, "=(E2=F2)")
$fc.Interior.Color = 5296274
